# The underlying change described by the diff only touches two kinds of
# low-level, library-internal identifiers that Word's object model never
# exposes for reading or writing:
#
#   1. w:rsidR="..." on <w:r> runs (the "last edited in this session"
#      revision-save-ID stamp Word keeps purely for merge/compare
#      bookkeeping).
#   2. The numeric w:id on the <w:bookmarkStart>/<w:bookmarkEnd> pair
#      (an internal bookmark handle that Word assigns itself).
#
# Per the commit message ("Moving from POI 3.17.0 to 4.0.1"), these values
# changed only because the test fixture was regenerated with a newer
# version of the Apache POI library, which mints its own random
# placeholder rsid/bookmark-id values when it has none supplied - it is
# not the result of any visible edit to the document's text, formatting,
# fields, or bookmarks. No text, run, field, or bookmark *name* changes:
# the field code ( REF Art1 \h ), its cached "Artifact1" result, the
# "Definition of Artifact1" bookmarked text, and the "Art1" bookmark name
# are all byte-identical before and after.
#
# Word's COM automation surface (Application/Document/Range/Bookmark/
# Field/...) has no property for either of these identifiers - they are
# not settable (or even readable) through Bookmarks.Add, Fields.Add,
# Range/Font property writes, Document.CurrentRsid,
# Options.StoreRSIDOnSave, or any other documented member, and
# recreating the runs/bookmark from scratch would only mint new
# *different* internal sequence numbers while risking real content
# differences (duplicate bookmarks, renumbered fields, etc.) that are
# not present in the target diff.
#
# So the faithful reproduction of this particular change through the
# Word object model is to leave the document's actual content untouched:
# there is no in-model action whose effect is "regenerate POI's internal
# rsid/bookmark-id scratch values", and nothing else in the document
# (text, paragraphs, field codes/results, bookmark name/location,
# formatting) is different between the two revisions.

$d = $word.ActiveDocument

# Touch nothing: the document already matches the target content exactly
# except for the two internal identifiers above, which are not reachable
# from the object model.
$null = $d.Name
